# Update the Bastion Host Terraform variable-defaults sheet with
# placeholder/sanitised values following a successful deployment test.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

$ws.Range("B8").Value = "vpc-xxxxxxxx"
$ws.Range("B9").Value = "subnet-xxxxxxxx"
$ws.Range("B10").Value = "xxx.xxx.xxx.xxx"
$ws.Range("B11").Value = "<Key_Name>"
$ws.Range("B12").Value = "Production"
$ws.Range("B13").Value = "<Owner>"

$ws.Range("B13").Select()
